$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 693.36365
$ws.Range("I2").Value = 80.888885
$ws.Range("K2").Value = 80.888885
$ws.Range("M2").Value = 32.111115
$ws.Range("H12").Value = 150
$ws.Range("I12").Value = 150
$ws.Range("K12").Value = 150
$ws.Range("M12").Value = 20
$ws.Range("H15").Value = 956.82855
$ws.Range("I15").Value = 956.82855
$ws.Range("K15").Value = 2870.48565
$ws.Range("M15").Value = -2701.48565
$ws.Range("H17").Value = 2090.2778
$ws.Range("J17").Value = 2090.2778
$ws.Range("L17").Value = 6270.8334
$ws.Range("N17").Value = -6606.8334
$ws.Range("H21").Value = 10999.5
$ws.Range("I21").Value = 10999.5
$ws.Range("K21").Value = 10999.5
$ws.Range("M21").Value = -10531.5
$ws.Range("H23").Value = 10999.5
$ws.Range("I23").Value = 10999.5
$ws.Range("K23").Value = 10999.5
$ws.Range("M23").Value = -10765.5
$ws.Range("H40").Value = 2762.7693
$ws.Range("I40").Value = 2209.6667
$ws.Range("J40").Value = 3236.8572
$ws.Range("K40").Value = 2209.6667
$ws.Range("L40").Value = 3236.8572
$ws.Range("M40").Value = -2034.6667
$ws.Range("N40").Value = -3586.8572
$ws.Range("H43").Value = 1247.5
$ws.Range("I43").Value = 1245
$ws.Range("J43").Value = 1250
$ws.Range("K43").Value = 1245
$ws.Range("L43").Value = 1250
$ws.Range("M43").Value = -1176
$ws.Range("N43").Value = -1388
$ws.Range("H51").Value = 3285.1428
$ws.Range("J51").Value = 3001.1667
$ws.Range("L51").Value = 3001.1667
$ws.Range("N51").Value = -3969.1667
$ws.Range("H53").Value = 1696.1111
$ws.Range("I53").Value = 1377.1666
$ws.Range("K53").Value = 1377.1666
$ws.Range("M53").Value = -740.1666
$ws.Range("H54").Value = 7002.9414
$ws.Range("I54").Value = 9005.556
$ws.Range("J54").Value = 4750
$ws.Range("K54").Value = 9005.556
$ws.Range("L54").Value = 4750
$ws.Range("M54").Value = -8519.556
$ws.Range("N54").Value = -5722
$ws.Range("H55").Value = 461.85715
$ws.Range("I55").Value = 732.5
$ws.Range("K55").Value = 732.5
$ws.Range("M55").Value = -518.5
$ws.Range("H70").Value = 3910.5557
$ws.Range("J70").Value = 4670.7144
$ws.Range("L70").Value = 14012.1432
$ws.Range("N70").Value = -14552.1432
$ws.Range("H73").Value = 3910.5557
$ws.Range("J73").Value = 4670.7144
$ws.Range("L73").Value = 14012.1432
$ws.Range("N73").Value = -15884.1432
$ws.Range("H76").Value = 4870.7827
$ws.Range("I76").Value = 4563
$ws.Range("J76").Value = 5978.8
$ws.Range("K76").Value = 4563
$ws.Range("L76").Value = 5978.8
$ws.Range("M76").Value = -4248
$ws.Range("N76").Value = -6608.8
$ws.Range("H79").Value = 4870.7827
$ws.Range("I79").Value = 4563
$ws.Range("J79").Value = 5978.8
$ws.Range("K79").Value = 4563
$ws.Range("L79").Value = 5978.8
$ws.Range("M79").Value = -3471
$ws.Range("N79").Value = -8162.8
$ws.Range("H86").Value = 3037
$ws.Range("J86").Value = 3025.25
$ws.Range("L86").Value = 3025.25
$ws.Range("N86").Value = -5271.25
$ws.Range("H89").Value = 3037
$ws.Range("J89").Value = 3025.25
$ws.Range("L89").Value = 15126.25
$ws.Range("N89").Value = -26358.25
$ws.Range("H94").Value = 17884.75
$ws.Range("I94").Value = 775
$ws.Range("K94").Value = 775
$ws.Range("M94").Value = -324
$ws.Range("H98").Value = 1649.3334
$ws.Range("I98").Value = 574.125
$ws.Range("J98").Value = 3799.75
$ws.Range("K98").Value = 574.125
$ws.Range("L98").Value = 3799.75
$ws.Range("M98").Value = 923.875
$ws.Range("N98").Value = -6795.75
$ws.Range("H100").Value = 978.6
$ws.Range("I100").Value = 973.5
$ws.Range("K100").Value = 973.5
$ws.Range("M100").Value = -432.5
$ws.Range("H113").Value = 5576.25
$ws.Range("I113").Value = 4646.95
$ws.Range("J113").Value = 7899.5
$ws.Range("K113").Value = 4646.95
$ws.Range("L113").Value = 7899.5
$ws.Range("M113").Value = -1392.95
$ws.Range("N113").Value = -14407.5
$ws.Range("H116").Value = 5677.4287
$ws.Range("I116").Value = 4957
$ws.Range("K116").Value = 4957
$ws.Range("M116").Value = -1515
$ws.Range("H121").Value = 3999
$ws.Range("J121").Value = 3999
$ws.Range("L121").Value = 11997
$ws.Range("N121").Value = -15491
$ws.Range("H122").Value = 1649.3334
$ws.Range("I122").Value = 574.125
$ws.Range("J122").Value = 3799.75
$ws.Range("K122").Value = 1722.375
$ws.Range("L122").Value = 11399.25
$ws.Range("M122").Value = 727.625
$ws.Range("N122").Value = -16299.25
$ws.Range("H125").Value = 1973.6
$ws.Range("I125").Value = 3444.5
$ws.Range("J125").Value = 993
$ws.Range("K125").Value = 31000.5
$ws.Range("L125").Value = 8937
$ws.Range("M125").Value = -28540.5
$ws.Range("N125").Value = -13857
$ws.Range("H135").Value = 11770
$ws.Range("I135").Value = 2025.125
$ws.Range("J135").Value = 50749.5
$ws.Range("K135").Value = 18226.125
$ws.Range("L135").Value = 456745.5
$ws.Range("M135").Value = -15691.125
$ws.Range("N135").Value = -461815.5
$ws.Range("H137").Value = 2937.3784
$ws.Range("I137").Value = 2226.2778
$ws.Range("K137").Value = 6678.8334
$ws.Range("M137").Value = -4128.8334
$ws.Range("H138").Value = 2273.2
$ws.Range("I138").Value = 1507.0646
$ws.Range("J138").Value = 2617.4058
$ws.Range("K138").Value = 4521.1938
$ws.Range("L138").Value = 7852.2174
$ws.Range("M138").Value = 618.8062
$ws.Range("N138").Value = -18132.2174
$ws.Range("H141").Value = 3521.75
$ws.Range("I141").Value = 2202.2354
$ws.Range("K141").Value = 6606.706200000001
$ws.Range("M141").Value = -1426.706200000001

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 663.3333
$ws.Range("I5").Value = 663.3333
$ws.Range("K5").Value = 663.3333
$ws.Range("M5").Value = -551.3333
$ws.Range("H32").Value = 2157.6052
$ws.Range("I32").Value = 1164.0422
$ws.Range("J32").Value = 16266.2
$ws.Range("K32").Value = 1164.0422
$ws.Range("L32").Value = 16266.2
$ws.Range("M32").Value = -877.0422000000001
$ws.Range("N32").Value = -16840.2
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 4000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H45").Value = 1937.8
$ws.Range("I45").Value = 2250
$ws.Range("K45").Value = 2250
$ws.Range("M45").Value = -1873
$ws.Range("H53").Value = 35009.75
$ws.Range("J53").Value = 39333.332
$ws.Range("L53").Value = 39333.332
$ws.Range("N53").Value = -40697.332
$ws.Range("H61").Value = 3222.4546
$ws.Range("I61").Value = 3109.3809
$ws.Range("J61").Value = 5597
$ws.Range("K61").Value = 3109.3809
$ws.Range("L61").Value = 5597
$ws.Range("M61").Value = -2897.3809
$ws.Range("N61").Value = -6021
$ws.Range("H74").Value = 1117.75
$ws.Range("I74").Value = 1141.8462
$ws.Range("K74").Value = 1141.8462
$ws.Range("M74").Value = -267.8462
$ws.Range("H77").Value = 1117.75
$ws.Range("I77").Value = 1141.8462
$ws.Range("K77").Value = 5709.231
$ws.Range("M77").Value = -1341.231
$ws.Range("H88").Value = 2235.8
$ws.Range("I88").Value = 2259.6667
$ws.Range("K88").Value = 2259.6667
$ws.Range("M88").Value = -1853.6667
$ws.Range("H91").Value = 2235.8
$ws.Range("I91").Value = 2259.6667
$ws.Range("K91").Value = 2259.6667
$ws.Range("M91").Value = -855.6667000000002
$ws.Range("H104").Value = 99997
$ws.Range("J104").Value = 99997
$ws.Range("L104").Value = 99997
$ws.Range("N104").Value = -106985
$ws.Range("H110").Value = 1601.7142
$ws.Range("I110").Value = 1590.4706
$ws.Range("J110").Value = 1649.5
$ws.Range("K110").Value = 1590.4706
$ws.Range("L110").Value = 1649.5
$ws.Range("M110").Value = 454.5293999999999
$ws.Range("N110").Value = -5739.5
$ws.Range("H122").Value = 3058.1924
$ws.Range("I122").Value = 3072.3809
$ws.Range("J122").Value = 2998.6
$ws.Range("K122").Value = 9217.1427
$ws.Range("L122").Value = 8995.799999999999
$ws.Range("M122").Value = -6767.1427
$ws.Range("N122").Value = -13895.8
$ws.Range("H132").Value = 1491.1794
$ws.Range("I132").Value = 1425.1578
$ws.Range("K132").Value = 4275.4734
$ws.Range("M132").Value = -1745.4734
$ws.Range("H136").Value = 3222.4546
$ws.Range("I136").Value = 3109.3809
$ws.Range("J136").Value = 5597
$ws.Range("K136").Value = 9328.1427
$ws.Range("L136").Value = 16791
$ws.Range("M136").Value = -6778.1427
$ws.Range("N136").Value = -21891

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 663.3333
$ws.Range("I4").Value = 663.3333
$ws.Range("K4").Value = 663.3333
$ws.Range("M4").Value = -548.3333
$ws.Range("H20").Value = 4633.25
$ws.Range("I20").Value = 2924.8462
$ws.Range("J20").Value = 12036.333
$ws.Range("K20").Value = 2924.8462
$ws.Range("L20").Value = 12036.333
$ws.Range("M20").Value = -2677.8462
$ws.Range("N20").Value = -12530.333
$ws.Range("H30").Value = 9980
$ws.Range("J30").Value = 9980
$ws.Range("L30").Value = 9980
$ws.Range("N30").Value = -10230
$ws.Range("H86").Value = 3575.3333
$ws.Range("I86").Value = 2796.25
$ws.Range("K86").Value = 2796.25
$ws.Range("M86").Value = -1673.25
$ws.Range("H89").Value = 3575.3333
$ws.Range("I89").Value = 2796.25
$ws.Range("K89").Value = 13981.25
$ws.Range("M89").Value = -8365.25
$ws.Range("H107").Value = 1204.3684
$ws.Range("I107").Value = 882.3889
$ws.Range("K107").Value = 882.3889
$ws.Range("M107").Value = 1037.6111
$ws.Range("H134").Value = 966.2593000000001
$ws.Range("I134").Value = 739.5599999999999
$ws.Range("K134").Value = 2218.68
$ws.Range("M134").Value = 316.3200000000002

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 77077.62
$ws.Range("I7").Value = 111245.336
$ws.Range("K7").Value = 111245.336
$ws.Range("M7").Value = -111132.336
$ws.Range("H16").Value = 887
$ws.Range("I16").Value = 774
$ws.Range("K16").Value = 774
$ws.Range("M16").Value = -487
$ws.Range("H31").Value = 1826.2142
$ws.Range("I31").Value = 1881
$ws.Range("J31").Value = 1497.5
$ws.Range("K31").Value = 1881
$ws.Range("L31").Value = 1497.5
$ws.Range("M31").Value = -1586
$ws.Range("N31").Value = -2087.5
$ws.Range("H34").Value = 1826.2142
$ws.Range("I34").Value = 1881
$ws.Range("J34").Value = 1497.5
$ws.Range("K34").Value = 1881
$ws.Range("L34").Value = 1497.5
$ws.Range("M34").Value = -1679
$ws.Range("N34").Value = -1901.5
$ws.Range("H55").Value = 35000
$ws.Range("I55").Value = 35000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 35000
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("M55").Value = -34685
$ws.Range("H62").Value = 7599.143
$ws.Range("I62").Value = 7299.6
$ws.Range("J62").Value = 8348
$ws.Range("K62").Value = 7299.6
$ws.Range("L62").Value = 8348
$ws.Range("M62").Value = -6675.6
$ws.Range("N62").Value = -9596
$ws.Range("H65").Value = 7599.143
$ws.Range("I65").Value = 7299.6
$ws.Range("J65").Value = 8348
$ws.Range("K65").Value = 36498
$ws.Range("L65").Value = 41740
$ws.Range("M65").Value = -33378
$ws.Range("N65").Value = -47980
$ws.Range("H99").Value = 2547
$ws.Range("I99").Value = 2547
$ws.Range("K99").Value = 2547
$ws.Range("M99").Value = -1049
$ws.Range("H113").Value = 887
$ws.Range("I113").Value = 774
$ws.Range("K113").Value = 774
$ws.Range("M113").Value = 1396
$ws.Range("H122").Value = 3245.625
$ws.Range("I122").Value = 1798.5
$ws.Range("K122").Value = 5395.5
$ws.Range("M122").Value = -2945.5
$ws.Range("H126").Value = 2547
$ws.Range("I126").Value = 2547
$ws.Range("K126").Value = 7641
$ws.Range("M126").Value = -5171
$ws.Range("H132").Value = 1030.5
$ws.Range("I132").Value = 671
$ws.Range("K132").Value = 2013
$ws.Range("M132").Value = 517

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 334.85715
$ws.Range("I8").Value = 334.85715
$ws.Range("K8").Value = 1004.57145
$ws.Range("M8").Value = -865.5714499999999
$ws.Range("H12").Value = 252.38095
$ws.Range("J12").Value = 160.45454
$ws.Range("L12").Value = 481.36362
$ws.Range("N12").Value = -827.3636200000001
$ws.Range("H14").Value = 412.7143
$ws.Range("I14").Value = 412.7143
$ws.Range("K14").Value = 1238.1429
$ws.Range("M14").Value = -1065.1429
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H23").Value = 269
$ws.Range("J23").Value = 290.9
$ws.Range("L23").Value = 872.6999999999999
$ws.Range("N23").Value = -1342.7
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H33").Value = 738
$ws.Range("I33").Value = 738
$ws.Range("K33").Value = 4428
$ws.Range("M33").Value = -4145
$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 3000
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118
$ws.Range("H55").Value = 2222.2222
$ws.Range("J55").Value = 2375
$ws.Range("L55").Value = 7125
$ws.Range("N55").Value = -7479
$ws.Range("H68").Value = 2527
$ws.Range("J68").Value = 750
$ws.Range("L68").Value = 2250
$ws.Range("N68").Value = -3872
$ws.Range("H71").Value = 2527
$ws.Range("J71").Value = 750
$ws.Range("L71").Value = 6750
$ws.Range("N71").Value = -14862
$ws.Range("H82").Value = 5664.8335
$ws.Range("I82").Value = 744.5
$ws.Range("J82").Value = 8125
$ws.Range("K82").Value = 2233.5
$ws.Range("L82").Value = 24375
$ws.Range("M82").Value = -1827.5
$ws.Range("N82").Value = -25187
$ws.Range("H85").Value = 5664.8335
$ws.Range("I85").Value = 744.5
$ws.Range("J85").Value = 8125
$ws.Range("K85").Value = 2233.5
$ws.Range("L85").Value = 24375
$ws.Range("M85").Value = -829.5
$ws.Range("N85").Value = -27183
$ws.Range("H113").Value = 515.3158
$ws.Range("I113").Value = 471.75
$ws.Range("J113").Value = 547
$ws.Range("K113").Value = 1415.25
$ws.Range("L113").Value = 1641
$ws.Range("M113").Value = 754.75
$ws.Range("N113").Value = -5981
$ws.Range("H114").Value = 2886.8
$ws.Range("I114").Value = 2727.25
$ws.Range("J114").Value = 2993.1667
$ws.Range("K114").Value = 8181.75
$ws.Range("L114").Value = 8979.500100000001
$ws.Range("M114").Value = -4927.75
$ws.Range("N114").Value = -15487.5001
$ws.Range("H117").Value = 2590.4285
$ws.Range("I117").Value = 598
$ws.Range("J117").Value = 2922.5
$ws.Range("K117").Value = 1794
$ws.Range("L117").Value = 8767.5
$ws.Range("M117").Value = 1648
$ws.Range("N117").Value = -15651.5
$ws.Range("H128").Value = 143661.42
$ws.Range("I128").Value = 143661.42
$ws.Range("K128").Value = 430984.26
$ws.Range("M128").Value = -426004.26
$ws.Range("H129").Value = 2096.6206
$ws.Range("I129").Value = 703.0909
$ws.Range("J129").Value = 2948.2222
$ws.Range("K129").Value = 2109.2727
$ws.Range("L129").Value = 8844.6666
$ws.Range("M129").Value = 2890.7273
$ws.Range("N129").Value = -18844.6666
$ws.Range("H130").Value = 4999.5
$ws.Range("J130").Value = 4999.5
$ws.Range("L130").Value = 14998.5
$ws.Range("N130").Value = -25038.5
$ws.Range("H131").Value = 28985.61
$ws.Range("I131").Value = 124025.89
$ws.Range("J131").Value = 2255.5312
$ws.Range("K131").Value = 372077.67
$ws.Range("L131").Value = 6766.5936
$ws.Range("M131").Value = -367037.67
$ws.Range("N131").Value = -16846.5936
$ws.Range("H140").Value = 1307.4348
$ws.Range("I140").Value = 928.55
$ws.Range("K140").Value = 2785.65
$ws.Range("M140").Value = 2394.35
$ws.Range("H141").Value = 2427.2856
$ws.Range("J141").Value = 2748.25
$ws.Range("L141").Value = 8244.75
$ws.Range("N141").Value = -18604.75

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 23000
$ws.Range("I55").Value = 23000
$ws.Range("K55").Value = 23000
$ws.Range("M55").Value = -22673
$ws.Range("H70").Value = 8404.454
$ws.Range("I70").Value = 5900
$ws.Range("J70").Value = 8654.9
$ws.Range("K70").Value = 5900
$ws.Range("L70").Value = 8654.9
$ws.Range("M70").Value = -5630
$ws.Range("N70").Value = -9194.9
$ws.Range("H73").Value = 8404.454
$ws.Range("I73").Value = 5900
$ws.Range("J73").Value = 8654.9
$ws.Range("K73").Value = 5900
$ws.Range("L73").Value = 8654.9
$ws.Range("M73").Value = -4964
$ws.Range("N73").Value = -10526.9
$ws.Range("H80").Value = 13521.45
$ws.Range("I80").Value = 6068.143
$ws.Range("J80").Value = 17534.77
$ws.Range("K80").Value = 6068.143
$ws.Range("L80").Value = 17534.77
$ws.Range("M80").Value = -5070.143
$ws.Range("N80").Value = -19530.77
$ws.Range("H83").Value = 13521.45
$ws.Range("I83").Value = 6068.143
$ws.Range("J83").Value = 17534.77
$ws.Range("K83").Value = 30340.715
$ws.Range("L83").Value = 87673.85000000001
$ws.Range("M83").Value = -25348.715
$ws.Range("N83").Value = -97657.85000000001
$ws.Range("H102").Value = 3417.25
$ws.Range("I102").Value = 2226.5
$ws.Range("K102").Value = 2226.5
$ws.Range("M102").Value = -604.5
$ws.Range("H122").Value = 2173.8262
$ws.Range("I122").Value = 2000.0526
$ws.Range("J122").Value = 2999.25
$ws.Range("K122").Value = 6000.1578
$ws.Range("L122").Value = 8997.75
$ws.Range("M122").Value = -3550.1578
$ws.Range("N122").Value = -13897.75
$ws.Range("H132").Value = 1603.5
$ws.Range("I132").Value = 947
$ws.Range("J132").Value = 4098.2
$ws.Range("K132").Value = 2841
$ws.Range("L132").Value = 12294.6
$ws.Range("M132").Value = -311
$ws.Range("N132").Value = -17354.6

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3362.5454
$ws.Range("J7").Value = 3873.75
$ws.Range("L7").Value = 3873.75
$ws.Range("N7").Value = -4097.75
$ws.Range("H16").Value = 409.30435
$ws.Range("I16").Value = 358.64706
$ws.Range("J16").Value = 552.8333
$ws.Range("K16").Value = 358.64706
$ws.Range("L16").Value = 552.8333
$ws.Range("M16").Value = -188.64706
$ws.Range("N16").Value = -892.8333
$ws.Range("H22").Value = 1287
$ws.Range("I22").Value = 1287
$ws.Range("K22").Value = 1287
$ws.Range("M22").Value = -992
$ws.Range("H27").Value = 1287
$ws.Range("I27").Value = 1287
$ws.Range("K27").Value = 1287
$ws.Range("M27").Value = -1180
$ws.Range("H32").Value = 7166.3335
$ws.Range("I32").Value = 5750
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 5750
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -5433
$ws.Range("N32").Value = -10633
$ws.Range("H40").Value = 3190.4285
$ws.Range("I40").Value = 2706.6
$ws.Range("J40").Value = 4400
$ws.Range("K40").Value = 2706.6
$ws.Range("L40").Value = 4400
$ws.Range("M40").Value = -2570.6
$ws.Range("N40").Value = -4672
$ws.Range("H43").Value = 4625240.5
$ws.Range("I43").Value = 2106666.8
$ws.Range("J43").Value = 6950078
$ws.Range("K43").Value = 2106666.8
$ws.Range("L43").Value = 6950078
$ws.Range("M43").Value = -2106473.8
$ws.Range("N43").Value = -6950464
$ws.Range("H53").Value = 40000
$ws.Range("J53").Value = 40000
$ws.Range("L53").Value = 40000
$ws.Range("N53").Value = -41036
$ws.Range("H55").Value = 587.9231
$ws.Range("I55").Value = 556.4167
$ws.Range("K55").Value = 556.4167
$ws.Range("M55").Value = -383.4167
$ws.Range("H68").Value = 2967.4375
$ws.Range("I68").Value = 2884.7144
$ws.Range("J68").Value = 3031.7778
$ws.Range("K68").Value = 2884.7144
$ws.Range("L68").Value = 3031.7778
$ws.Range("M68").Value = -2135.7144
$ws.Range("N68").Value = -4529.7778
$ws.Range("H71").Value = 2967.4375
$ws.Range("I71").Value = 2884.7144
$ws.Range("J71").Value = 3031.7778
$ws.Range("K71").Value = 14423.572
$ws.Range("L71").Value = 15158.889
$ws.Range("M71").Value = -10679.572
$ws.Range("N71").Value = -22646.889
$ws.Range("H93").Value = 3500.5
$ws.Range("I93").Value = 3431.3076
$ws.Range("J93").Value = 4400
$ws.Range("K93").Value = 3431.3076
$ws.Range("L93").Value = 4400
$ws.Range("M93").Value = -2183.3076
$ws.Range("N93").Value = -6896
$ws.Range("H108").Value = 75184.60000000001
$ws.Range("J108").Value = 75184.60000000001
$ws.Range("L108").Value = 75184.60000000001
$ws.Range("N108").Value = -82864.60000000001
$ws.Range("H126").Value = 3362.5454
$ws.Range("J126").Value = 3873.75
$ws.Range("L126").Value = 11621.25
$ws.Range("N126").Value = -16561.25
$ws.Range("H132").Value = 1913.4
$ws.Range("I132").Value = 1912.1428
$ws.Range("J132").Value = 1920
$ws.Range("K132").Value = 5736.428400000001
$ws.Range("L132").Value = 5760
$ws.Range("M132").Value = -3206.428400000001
$ws.Range("N132").Value = -10820
$ws.Range("H136").Value = 2512.0908
$ws.Range("I136").Value = 2614.3333
$ws.Range("J136").Value = 2052
$ws.Range("K136").Value = 7842.999899999999
$ws.Range("L136").Value = 6156
$ws.Range("M136").Value = -5292.999899999999
$ws.Range("N136").Value = -11256

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 5101
$ws.Range("I9").Value = 5101
$ws.Range("K9").Value = 5101
$ws.Range("M9").Value = -4961
$ws.Range("H14").Value = 3745
$ws.Range("J14").Value = 4990
$ws.Range("L14").Value = 4990
$ws.Range("N14").Value = -5326
$ws.Range("H54").Value = 22721.666
$ws.Range("J54").Value = 18000
$ws.Range("L54").Value = 18000
$ws.Range("N54").Value = -19040
$ws.Range("H74").Value = 29750.334
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 29750.334
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 29750.334
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -31622.334
$ws.Range("H77").Value = 29750.334
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 29750.334
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 89251.00199999999
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -98611.00199999999
$ws.Range("H107").Value = 898.6
$ws.Range("I107").Value = 1021.1429
$ws.Range("J107").Value = 612.6667
$ws.Range("K107").Value = 3063.4287
$ws.Range("L107").Value = 1838.0001
$ws.Range("M107").Value = -1143.4287
$ws.Range("N107").Value = -5678.0001
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 1181
$ws.Range("I113").Value = 1296.4286
$ws.Range("J113").Value = 777
$ws.Range("K113").Value = 3889.2858
$ws.Range("L113").Value = 2331
$ws.Range("M113").Value = -1719.2858
$ws.Range("N113").Value = -6671
$ws.Range("H122").Value = 2741.75
$ws.Range("I122").Value = 2186.3215
$ws.Range("K122").Value = 6558.9645
$ws.Range("M122").Value = -4108.9645
$ws.Range("H132").Value = 1623
$ws.Range("I132").Value = 1623
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4869
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2339
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1419.2667
$ws.Range("J136").Value = 1985.3334
$ws.Range("L136").Value = 5956.0002
$ws.Range("N136").Value = -11056.0002

Write-Output "All changes applied"